$wb = $excel.ActiveWorkbook

# New data (rows 2-16, columns B and C) for sheet "NBR"
$nbrB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$nbrC = @(870,861,835,832,828,822,819,813,774,774,753,749,751,749,747)

# New data (rows 2-16, columns B and C) for sheet "BAR"
$barB = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
$barC = @(756,752,753,751,761,758,756,757,757,758,758,755,751,751,753)

foreach ($sheetInfo in @(
        @{ Name = "NBR"; B = $nbrB; C = $nbrC },
        @{ Name = "BAR"; B = $barB; C = $barC }
    )) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    # Remove the last 4 data rows (rows 17-20), shrinking the table from 19 to 15 data rows
    $ws.Rows("17:20").Delete()

    # Update the remaining data rows (2-16) with the new column B / C values
    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $sheetInfo.B[$i]
        $ws.Cells.Item($row, 3).Value = $sheetInfo.C[$i]
    }
}
